$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, C, E, F, G, I, K, M, O hold the per-row loading-percent series
# (D, H, J, L, N stay at their existing 0 values and are left untouched).
$cols = @("B", "C", "E", "F", "G", "I", "K", "M", "O")

# New values for rows 2..25 (one row per array), in the same column order as $cols.
$newValues = @(
    @(8.081187044361265, 4.868885135800753, 12.78731148306921, 16.86991607391245, 3.63622180581312, 20.65332281824074, 7.919350851870529, 13.37601256877296, 21.65688099826048),
    @(7.764384920425327, 4.705078792479989, 12.5827115674136, 15.89584955866815, 3.637859404668161, 20.77057059430543, 7.709206396593488, 13.19755499872676, 21.77027105084723),
    @(7.563795508589352, 4.600707676867058, 12.46007882931824, 15.26997757108489, 3.638917340391958, 20.84679552166092, 7.577983424331817, 13.08977309511277, 21.84485241634784),
    @(7.480645052716572, 4.557261801265673, 12.41092938550938, 15.008197319934, 3.639361686306258, 20.87892313009361, 7.524032995501652, 13.04635395027177, 21.87649025580648),
    @(7.466756623603518, 4.549993710492458, 12.40282001765264, 14.96433081551589, 3.639436269788074, 20.884322256456, 7.515048112217919, 13.03917611702196, 21.88181885938717),
    @(7.562679649933698, 4.600125394446204, 12.45941255100485, 15.26647399323133, 3.638923279379793, 20.84722449154397, 7.577257654247505, 13.08918542483904, 21.84527405533155),
    @(7.973280549737491, 4.813214753169545, 12.71618796220288, 16.5399640634477, 3.636775590109499, 20.69287135938985, 7.847392235604052, 13.31414013010204, 21.69494749513488),
    @(8.725770382474165, 5.199434799866355, 13.24024575756554, 19.00274580682531, 3.632978207025967, 20.42375634696059, 8.461369263491095, 13.76716444605257, 21.43958783509888),
    @(9.241308432242089, 5.462099524199733, 13.6334402047557, 20.67494806633232, 3.63043814703801, 20.24646654781685, 8.906496779368936, 14.10411944615155, 21.2761233282384),
    @(9.46686374841493, 5.576704370393052, 13.81317920606532, 21.3917225636224, 3.629336302730785, 20.17024224963027, 9.100380106200049, 14.25764597386327, 21.20702775455759),
    @(9.550926647506772, 5.619378554363141, 13.88129351382618, 21.65686569030329, 3.628926734013052, 20.14201428764682, 9.172526594130055, 14.31576385541119, 21.18162251488705),
    @(9.532883116730444, 5.610220434675051, 13.86662276917189, 21.60004134736742, 3.629014601179677, 20.14806537206263, 9.157045679666977, 14.30324883835198, 21.18706015509333),
    @(9.473807027430295, 5.580229831162534, 13.81878240628663, 21.4136618050453, 3.629302453646944, 20.1679071603719, 9.106341336788086, 14.26242805994999, 21.20492240834684),
    @(9.437443749900352, 5.561764804430171, 13.78948328023, 21.29868154950795, 3.629479770094389, 20.180143727693, 9.07511672936463, 14.23742005295441, 21.21596257557214),
    @(9.226382058339036, 5.45450962426413, 13.62170498536032, 20.62722412089977, 3.630511231324578, 20.25153704960099, 8.893649725188173, 14.09408635863009, 21.28074505897242),
    @(9.09456021868864, 5.38744485342557, 13.51894831568335, 20.20408069597325, 3.631157711317778, 20.29646812873607, 8.780094423485387, 14.00617918233049, 21.32183741180111),
    @(9.017898948708028, 5.348412517762764, 13.45993188560157, 19.95656407809801, 3.631534600764887, 20.32272785525238, 8.713972020324789, 13.95564344587492, 21.34596807557496),
    @(8.99180034129254, 5.335118819739844, 13.43996709651382, 19.87204792380568, 3.631663077770533, 20.33169048959024, 8.691446420626061, 13.93853908253213, 21.35422331556579),
    @(9.108680369486727, 5.394631629930313, 13.52987853404591, 20.24955283636154, 3.631088369907506, 20.29164202586517, 8.792266515594974, 14.01553474335382, 21.31741176685714),
    @(9.491196187421142, 5.589058614244775, 13.83283347695076, 21.46857628470577, 3.629217696380395, 20.16206186977399, 9.121269234557579, 14.27441906189158, 21.19965519085052),
    @(9.733301200778026, 5.711898512830331, 14.03109680066341, 22.22866616901552, 3.628039826233106, 20.08108401225912, 9.328858960998469, 14.44347834255036, 21.12712415581029),
    @(9.60482514584808, 5.646730099338121, 13.92527968472484, 21.82633154458858, 3.628664397942691, 20.12396387194684, 9.218754896877524, 14.35327831513417, 21.16542913397695),
    @(9.102299370472968, 5.391383968940066, 13.5249367905285, 20.22900810905287, 3.631119702897831, 20.29382257297934, 8.786766120088304, 14.0113050798615, 21.31941102456245),
    @(8.528424375226791, 5.098540067618897, 13.09672672491625, 18.34778573295695, 3.633961428487244, 20.49296915184413, 8.290077205034693, 13.64366052543064, 21.50443695813647)
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowValues = $newValues[$i]
    $rowNum = $startRow + $i
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range("$($cols[$j])$rowNum").Value = $rowValues[$j]
    }
}
